# Applies the "received tabular data for verification of state models" edit:
#  - updates clock-bias verification values (b1, b2) on initialConditions
#  - updates gravity/atmospheric acceleration error values on initialConditions
#  - updates position/velocity error-injection values on errorInjection
#  - updates the active sheet / selected cell to reflect where the user ended up

$wb = $excel.ActiveWorkbook

# --- initialConditions sheet -------------------------------------------------
$ic = $wb.Worksheets.Item("initialConditions")

# clock bias verification values (B8/B15 feed E8/E15 via "=Bn" formulas)
$ic.Range("B8").Value = 0.0000000547
$ic.Range("B15").Value = 0.00000000323

# gravity / atmospheric acceleration error values (B57:B62 feed E57:E62 via "=Bn")
$ic.Range("B57").Value = 0.01
$ic.Range("B58").Value = -0.01
$ic.Range("B59").Value = 0.001
$ic.Range("B60").Value = -0.07
$ic.Range("B61").Value = 0.1
$ic.Range("B62").Value = -0.04

# --- errorInjection sheet ----------------------------------------------------
$ei = $wb.Worksheets.Item("errorInjection")

# position / velocity error injection values (B2:B7 feed E2:E7 via "=Bn*1000")
$ei.Range("B2").Value = 25
$ei.Range("B3").Value = 47
$ei.Range("B4").Value = -32
$ei.Range("B5").Value = 16
$ei.Range("B6").Value = -29
$ei.Range("B7").Value = 13

# --- selection / active sheet -------------------------------------------------
# user finished editing on initialConditions (scrolled on to B63) so that
# sheet becomes the active tab; errorInjection (previously active) keeps its
# own last selection at B13
$ei.Range("B13").Select()

$ic.Activate()
$ic.Range("B63").Select()
